$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "Data" sheet: append 10 new daily observations (rows 463-472)
# ------------------------------------------------------------------
$dataWs = $wb.Worksheets.Item("Data")

$newRows = @(
    @(45147, 1796.519),
    @(45148, 1759.897),
    @(45149, 1773.236),
    @(45152, 1799.311),
    @(45153, 1743.784),
    @(45154, 1796.725),
    @(45155, 1794.12),
    @(45156, 1819.201),
    @(45159, 1824.788),
    @(45160, 1812.294)
)

$lastExistingRow = 462
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $dateSerial = $newRows[$i][0]
    $value = $newRows[$i][1]

    # Copy the date cell's formatting (matches the style used by the rest
    # of column A, a center-aligned YYYY-MM-DD HH:MM:SS date format) from
    # the preceding row, then set the actual values.
    $dataWs.Cells.Item($lastExistingRow, 1).Copy() | Out-Null
    $dataWs.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $dataWs.Cells.Item($r, 1).Value = $dateSerial
    $dataWs.Cells.Item($r, 2).Value = $value
}

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) "SeriesInfo" sheet: refresh the series metadata
# ------------------------------------------------------------------
$infoWs = $wb.Worksheets.Item("SeriesInfo")

function Set-PlainTextValue($range, $text) {
    # Writing a date-shaped string straight into a General cell makes
    # Excel silently convert it into a date serial number. Force the
    # cell to Text first so the literal string is kept, then drop the
    # formatting override again so the cell ends up unstyled, exactly
    # like it was before the edit.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-PlainTextValue $infoWs.Range("B3") "2023-08-22"
Set-PlainTextValue $infoWs.Range("B4") "2023-08-22"
Set-PlainTextValue $infoWs.Range("B7") "2023-08-22"

$infoWs.Range("B14").Value = "2023-08-22 13:01:06-05"
$infoWs.Range("B15").Value = 91
